$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormatLocal = "@"
$ws.Range("D2").Value = "66.295.85"
$ws.Range("E2").Value = "  +3.18%  "
$ws.Range("D3").NumberFormatLocal = "@"
$ws.Range("D3").Value = "3.490.48"
$ws.Range("E3").Value = "  +5.33%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormatLocal = "@"
$ws.Range("D5").Value = "559.06"
$ws.Range("E5").Value = "  +6.00%  "
$ws.Range("D6").NumberFormatLocal = "@"
$ws.Range("D6").Value = "182.67"
$ws.Range("E6").Value = "  +5.71%  "
$ws.Range("D7").NumberFormatLocal = "@"
$ws.Range("D7").Value = "0.639"
$ws.Range("E7").Value = "  +9.15%  "
$ws.Range("D8").NumberFormatLocal = "@"
$ws.Range("D8").Value = "3.491.20"
$ws.Range("E8").Value = "  +5.77%  "
$ws.Range("E10").Value = "  +5.31%  "
$ws.Range("D11").NumberFormatLocal = "@"
$ws.Range("D11").Value = "0.153"
$ws.Range("E11").Value = "  +14.83%  "
$ws.Range("D12").NumberFormatLocal = "@"
$ws.Range("D12").Value = "54.51"
$ws.Range("E12").Value = "  +3.02%  "
$ws.Range("D13").NumberFormatLocal = "@"
$ws.Range("D13").Value = "0.0000277"
$ws.Range("E13").Value = "  +7.80%  "
$ws.Range("D14").NumberFormatLocal = "@"
$ws.Range("D14").Value = "9.31"
$ws.Range("E14").Value = "  +4.31%  "
$ws.Range("D15").NumberFormatLocal = "@"
$ws.Range("D15").Value = "4.056.59"
$ws.Range("E15").Value = "  +5.88%  "
$ws.Range("D16").NumberFormatLocal = "@"
$ws.Range("D16").Value = "3.484.19"
$ws.Range("E16").Value = "  +5.32%  "
$ws.Range("D17").NumberFormatLocal = "@"
$ws.Range("D17").Value = "18.55"
$ws.Range("E17").Value = "  +6.51%  "
$ws.Range("E18").Value = "  +4.14%  "
$ws.Range("D19").NumberFormatLocal = "@"
$ws.Range("D19").Value = "66.282.29"
$ws.Range("E19").Value = "  +3.40%  "
$ws.Range("D20").NumberFormatLocal = "@"
$ws.Range("D20").Value = "12.02"
$ws.Range("E20").Value = "  +7.90%  "
$ws.Range("D21").NumberFormatLocal = "@"
$ws.Range("D21").Value = "0.996"
$ws.Range("E21").Value = "  +4.43%  "
$ws.Range("D22").NumberFormatLocal = "@"
$ws.Range("D22").Value = "418.84"
$ws.Range("E22").Value = "  +10.74%  "
$ws.Range("D23").NumberFormatLocal = "@"
$ws.Range("D23").Value = "4.07"
$ws.Range("E23").Value = "  +10.03%  "
$ws.Range("D24").NumberFormatLocal = "@"
$ws.Range("D24").Value = "86.15"
$ws.Range("D25").NumberFormatLocal = "@"
$ws.Range("D25").Value = "4.29"
$ws.Range("E25").Value = "  +3.04%  "
$ws.Range("E26").Value = "  +8.11%  "
$ws.Range("D27").NumberFormatLocal = "@"
$ws.Range("D27").Value = "10.81"
$ws.Range("E27").Value = "  -2.75%  "
$ws.Range("D28").NumberFormatLocal = "@"
$ws.Range("D28").Value = "12.34"
$ws.Range("E28").Value = "  +9.88%  "
$ws.Range("D29").NumberFormatLocal = "@"
$ws.Range("D29").Value = "6.04"
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("D30").NumberFormatLocal = "@"
$ws.Range("D30").Value = "9.07"
$ws.Range("E30").Value = "  +11.62%  "
$ws.Range("D31").NumberFormatLocal = "@"
$ws.Range("D31").Value = "30.26"
$ws.Range("E31").Value = "  +5.21%  "
$ws.Range("D32").NumberFormatLocal = "@"
$ws.Range("D32").Value = "6.81"
$ws.Range("E32").Value = "  +3.62%  "
$ws.Range("D33").NumberFormatLocal = "@"
$ws.Range("D33").Value = "622.81"
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("E34").Value = "  +5.71%  "
$ws.Range("E35").Value = "  +5.91%  "
$ws.Range("D36").NumberFormatLocal = "@"
$ws.Range("D36").Value = "60.31"
$ws.Range("E36").Value = "  +5.78%  "
$ws.Range("D37").NumberFormatLocal = "@"
$ws.Range("D37").Value = "0.149"
$ws.Range("E37").Value = "  +18.86%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").NumberFormatLocal = "@"
$ws.Range("D38").Value = "37.86"
$ws.Range("E38").Value = "  +5.27%  "
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").NumberFormatLocal = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").NumberFormatLocal = "@"
$ws.Range("D40").Value = "0.0₃0795"
$ws.Range("E40").Value = "  +6.09%  "
$ws.Range("E41").Value = "  +1.69%  "
$ws.Range("D42").NumberFormatLocal = "@"
$ws.Range("D42").Value = "3.38"
$ws.Range("E42").Value = "  +5.99%  "
$ws.Range("D43").NumberFormatLocal = "@"
$ws.Range("D43").Value = "3.109.23"
$ws.Range("E43").Value = "  +7.74%  "
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").NumberFormatLocal = "@"
$ws.Range("D45").Value = "2.85"
$ws.Range("E45").Value = "  +8.86%  "
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormatLocal = "@"
$ws.Range("D46").Value = "2.58"
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormatLocal = "@"
$ws.Range("D47").Value = "0.0416"
$ws.Range("E47").Value = "  +4.79%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormatLocal = "@"
$ws.Range("D48").Value = "3.27"
$ws.Range("E48").Value = "  +7.00%  "
$ws.Range("D49").NumberFormatLocal = "@"
$ws.Range("D49").Value = "2.75"
$ws.Range("E49").Value = "  +2.17%  "
$ws.Range("D50").NumberFormatLocal = "@"
$ws.Range("D50").Value = "0.133"
$ws.Range("D51").NumberFormatLocal = "@"
$ws.Range("D51").Value = "139.96"
$ws.Range("E51").Value = "  +1.57%  "
